$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 2 ---------------------------------
# This shifts the existing row 2 (game "IBJAWlHc") down to row 3, and
# leaves a blank row 2 ready to receive the new game's data.
$ws.Rows.Item(2).Insert()
# Row-insert in Excel copies formatting from the row above (the bold
# header row); strip that back off so the new data row has no style,
# matching the rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# --- 2. Add the two new trailing columns' headers -----------------------
$ws.Cells.Item(1, 55).Value = "Odd_CS_3-3_HT"
$ws.Cells.Item(1, 56).Value = "Odd_CS_4-4_HT"

# --- 3. Fill in the new row 2 (game "AiivoxjC") --------------------------
$row2 = @(
    "AiivoxjC", "25/11/2024", "12:30", "BULGARIA - PARVA LIGA", "Slavia Sofia", "Krumovgrad",
    1.8, 3.4, 4.5, 2.5, 2.05, 5, 1.08, 8, 1.4, 2.75, 2.3, 1.6, 1.5, 2.5, 2.1, 1.67,
    6, 7.5, 9, 15, 17, 34, 7.5, 6.5, 19, 67, 1250,
    10, 21, 15, 51, 41, 51,
    3.6, 10, 23, 34, 67, 201,
    2.5, 9, 67, 6, 26, 41, 101, 151, 351
)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# --- 4. Append the brand-new row 4 (game "j3x6GxT7") ----------------------
$row4 = @(
    "j3x6GxT7", "25/11/2024", "12:30", "ROMANIA - LIGA 1", "FC Botosani", "Poli Iasi",
    2.45, 2.8, 3.1, 3.25, 1.95, 3.75, 1.1, 7, 1.44, 2.63, 2.4, 1.53, 1.53, 2.38, 2, 1.73,
    7, 11, 10, 23, 23, 41, 7, 5.5, 15, 51, 351,
    8, 15, 12, 34, 29, 41,
    4.33, 15, 29, 51, 81, 251,
    2.38, 9, 67, 4.75, 19, 29, 51, 101, 251,
    51, 51
)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

Write-Host "edit complete"
